$p = $ppt.ActivePresentation
try {
    $p.ApplyTheme("doesnotexist.xml")
    Write-Host "Applied theme"
} catch {
    Write-Host "Error: $_"
}
